$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.225.24"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "2.240.32"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.95"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -3.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.26"
$ws.Range("E10").Value = "  -2.91%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.93"
$ws.Range("E12").Value = "  -2.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "2.573.76"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.36"
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("E16").Value = "  -1.98%  "

$ws.Range("D17").Value = "2.250.37"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "42.088.81"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("E19").Value = "  -5.51%  "

$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.73"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  +10.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.19"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("E24").Value = "  -6.09%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.40"
$ws.Range("E26").Value = "  -3.35%  "

$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("E29").Value = "  -2.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.61"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.65"
$ws.Range("E32").Value = "  -3.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0808"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.75"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").Value = "  -6.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("E37").Value = "  -5.04%  "

$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.16"
$ws.Range("E39").Value = "  -1.94%  "

$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.56"
$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.199"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.73"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.69"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").Value = "2.448.23"
$ws.Range("E51").Value = "  -0.06%  "
